# Gallery_Loops_P_Panels.xlsx - add Netherlands / Austria / Denmark market test-data sheets
$wb = $excel.ActiveWorkbook

# --- Netherlands (copied from the Belgium-shaped 10-row template) ---
$src = $wb.Worksheets.Item("Belgium")
$src.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$nl = $wb.Worksheets.Item($wb.Worksheets.Count)
$nl.Name = "Netherlands"
$nl.Range("B4").Value = "NGC-3144/T2188"
$nl.Range("B2").Value = "Netherlands Market"
$nl.Range("B4").Select() | Out-Null

# --- Austria (copied from the Germany-shaped 12-row template) ---
$src = $wb.Worksheets.Item("Germany")
$src.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$at = $wb.Worksheets.Item($wb.Worksheets.Count)
$at.Name = "Austria"
$at.Range("B4").Value = "NGC-3817/T2295"
$at.Range("B2").Value = "Netherlands Market"
$at.Range("D21").Select() | Out-Null

# --- Denmark (copied from the Belgium-shaped 10-row template) ---
$src = $wb.Worksheets.Item("Belgium")
$src.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$dk = $wb.Worksheets.Item($wb.Worksheets.Count)
$dk.Name = "Denmark"
$dk.Range("B4").Value = "NGC-2913/T2783"
$dk.Range("B2").Value = "Denmark Market"
$dk.Range("D17").Select() | Out-Null

# Denmark ends up the active tab; Greece's earlier selection moves to B14
$dk.Activate() | Out-Null
$wb.Worksheets.Item("Greece").Range("B14").Select() | Out-Null
$dk.Activate() | Out-Null
